$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these as text, preserving exact
# formatting (trailing zeros, multi-dot strings, no scientific notation).

$ws.Range("D2").Value = "`'26.225.86"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "`'1.671.91"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("D4").Value = "`'1.007"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").Value = "`'217.81"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").Value = "`'0.5124"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").Value = "`'1.007"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "`'0.2660"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").Value = "`'0.06381"
$ws.Range("E9").Value = "  +2.56%  "
$ws.Range("D10").Value = "`'21.54"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("D11").Value = "`'0.07396"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "`'1.683.23"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "`'4.550"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "`'0.5833"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "`'1.900.28"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "`'0.000008652"
$ws.Range("E16").Value = "  +4.65%  "
$ws.Range("D17").Value = "`'64.56"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "`'26.324.33"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "`'4.960"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "`'1.007"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "`'10.89"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").Value = "`'189.63"
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("D23").Value = "`'6.213"
$ws.Range("D24").Value = "`'1.008"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "`'144.07"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "`'7.660"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "`'0.1182"
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("E29").Value = "  +1.63%  "
$ws.Range("D30").Value = "`'1.282"
$ws.Range("E30").Value = "  -4.25%  "
$ws.Range("D31").Value = "`'1.327"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").Value = "`'3.529"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("D33").Value = "`'3.532"
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").Value = "`'0.6036"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "`'2.377"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "`'2.647"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "`'0.01618"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").Value = "`'6.093"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D41").Value = "`'1.081.94"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "`'0.8720"
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").Value = "`'1.011"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "`'100.43"
$ws.Range("E44").Value = "  +2.98%  "
$ws.Range("D45").Value = "`'1.823.71"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  +7.37%  "
$ws.Range("D47").Value = "`'56.34"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").Value = "`'8.097"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("D50").Value = "`'0.05220"
$ws.Range("D51").Value = "`'0.4299"
$ws.Range("E51").Value = "  -1.84%  "
